$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1.129247026424488
$ws.Cells.Item(2, 4).Value = 1.130722603127263
$ws.Cells.Item(2, 5).Value = 1.127063990637527
$ws.Cells.Item(2, 6).Value = 1.14098192733876
$ws.Cells.Item(2, 10).Value = 1.133896276542331
$ws.Cells.Item(2, 11).Value = 1.13327088222099
$ws.Cells.Item(2, 12).Value = 1.1296209710029
$ws.Cells.Item(2, 13).Value = 1.14350612594359
$ws.Cells.Item(2, 14).Value = 1.135506539314743
$ws.Cells.Item(3, 3).Value = 1.133225707344481
$ws.Cells.Item(3, 4).Value = 1.134456356757958
$ws.Cells.Item(3, 5).Value = 1.130691390990533
$ws.Cells.Item(3, 6).Value = 1.144817092886283
$ws.Cells.Item(3, 10).Value = 1.137541038352138
$ws.Cells.Item(3, 11).Value = 1.136825465177155
$ws.Cells.Item(3, 12).Value = 1.13306876888923
$ws.Cells.Item(3, 13).Value = 1.147163744512365
$ws.Cells.Item(3, 14).Value = 1.139156477104379
$ws.Cells.Item(4, 3).Value = 1.135778047712894
$ws.Cells.Item(4, 4).Value = 1.136850916933994
$ws.Cells.Item(4, 5).Value = 1.133017400812488
$ws.Cells.Item(4, 6).Value = 1.147277086463468
$ws.Cells.Item(4, 10).Value = 1.139877754619945
$ws.Cells.Item(4, 11).Value = 1.139104018246755
$ws.Cells.Item(4, 12).Value = 1.135278475611613
$ws.Cells.Item(4, 13).Value = 1.149508788926253
$ws.Cells.Item(4, 14).Value = 1.141496511777311
$ws.Cells.Item(5, 3).Value = 1.136845932583764
$ws.Cells.Item(5, 4).Value = 1.137852628575423
$ws.Cells.Item(5, 5).Value = 1.133990353537638
$ws.Cells.Item(5, 6).Value = 1.14830626510601
$ws.Cells.Item(5, 10).Value = 1.140855084841424
$ws.Cells.Item(5, 11).Value = 1.140056938218034
$ws.Cells.Item(5, 12).Value = 1.136202508391054
$ws.Cells.Item(5, 13).Value = 1.150489622383769
$ws.Cells.Item(5, 14).Value = 1.142475229919806
$ws.Cells.Item(6, 3).Value = 1.137024940099218
$ws.Cells.Item(6, 4).Value = 1.138020534239541
$ws.Cells.Item(6, 5).Value = 1.13415343382701
$ws.Cells.Item(6, 6).Value = 1.148478780285695
$ws.Cells.Item(6, 10).Value = 1.141018892970114
$ws.Cells.Item(6, 11).Value = 1.140216650092962
$ws.Cells.Item(6, 12).Value = 1.136357373131804
$ws.Cells.Item(6, 13).Value = 1.150654018841646
$ws.Cells.Item(6, 14).Value = 1.142639270674829
$ws.Cells.Item(7, 3).Value = 1.135792336721849
$ws.Cells.Item(7, 4).Value = 1.136864321128296
$ws.Cells.Item(7, 5).Value = 1.133030420500037
$ws.Cells.Item(7, 6).Value = 1.147290857830433
$ws.Cells.Item(7, 10).Value = 1.139890833281393
$ws.Cells.Item(7, 11).Value = 1.139116770578635
$ws.Cells.Item(7, 12).Value = 1.135290841737633
$ws.Cells.Item(7, 13).Value = 1.14952191439
$ws.Cells.Item(7, 14).Value = 1.141509609011959
$ws.Cells.Item(8, 3).Value = 1.13059633255364
$ws.Cells.Item(8, 4).Value = 1.131988982407891
$ws.Cells.Item(8, 5).Value = 1.12829436932283
$ws.Cells.Item(8, 6).Value = 1.142282620110824
$ws.Cells.Item(8, 10).Value = 1.135132633570559
$ws.Cells.Item(8, 11).Value = 1.134476721655373
$ws.Cells.Item(8, 12).Value = 1.130790667433169
$ws.Cells.Item(8, 13).Value = 1.14474682607886
$ws.Cells.Item(8, 14).Value = 1.136744652111762
$ws.Cells.Item(9, 3).Value = 1.121262506868623
$ws.Cells.Item(9, 4).Value = 1.123226112099647
$ws.Cells.Item(9, 5).Value = 1.119779189702397
$ws.Cells.Item(9, 6).Value = 1.133283969947175
$ws.Cells.Item(9, 10).Value = 1.126574254515697
$ws.Cells.Item(9, 11).Value = 1.126128148900321
$ws.Cells.Item(9, 12).Value = 1.122690686595263
$ws.Cells.Item(9, 13).Value = 1.136158758830098
$ws.Cells.Item(9, 14).Value = 1.128174119176983
$ws.Cells.Item(10, 3).Value = 1.114908459795771
$ws.Cells.Item(10, 4).Value = 1.117257403583155
$ws.Cells.Item(10, 5).Value = 1.113977409134178
$ws.Cells.Item(10, 6).Value = 1.127156767337604
$ws.Cells.Item(10, 10).Value = 1.120740681761555
$ws.Cells.Item(10, 11).Value = 1.12043576742126
$ws.Cells.Item(10, 12).Value = 1.1171657353171
$ws.Cells.Item(10, 13).Value = 1.130305446719618
$ws.Cells.Item(10, 14).Value = 1.12233226208041
$ws.Cells.Item(11, 3).Value = 1.112123128571427
$ws.Cells.Item(11, 4).Value = 1.114640213447464
$ws.Cells.Item(11, 5).Value = 1.111432994774486
$ws.Cells.Item(11, 6).Value = 1.124470582955499
$ws.Cells.Item(11, 10).Value = 1.118181747421183
$ws.Cells.Item(11, 11).Value = 1.117938338832573
$ws.Cells.Item(11, 12).Value = 1.114741274509416
$ws.Cells.Item(11, 13).Value = 1.127737978669283
$ws.Cells.Item(11, 14).Value = 1.11976969375976
$ws.Cells.Item(12, 3).Value = 1.111083180531546
$ws.Cells.Item(12, 4).Value = 1.113662928943057
$ws.Cells.Item(12, 5).Value = 1.110482822726678
$ws.Cells.Item(12, 6).Value = 1.123467611216252
$ws.Cells.Item(12, 10).Value = 1.117226064858702
$ws.Cells.Item(12, 11).Value = 1.117005562577925
$ws.Cells.Item(12, 12).Value = 1.113835678015828
$ws.Cells.Item(12, 13).Value = 1.126779128449707
$ws.Cells.Item(12, 14).Value = 1.118812654018422
$ws.Cells.Item(13, 3).Value = 1.111306499163178
$ws.Cells.Item(13, 4).Value = 1.113872796351518
$ws.Cells.Item(13, 5).Value = 1.110686870689613
$ws.Cells.Item(13, 6).Value = 1.123682991404905
$ws.Cells.Item(13, 10).Value = 1.117431300211459
$ws.Cells.Item(13, 11).Value = 1.117205881638884
$ws.Cells.Item(13, 12).Value = 1.114030163379014
$ws.Cells.Item(13, 13).Value = 1.126985043176907
$ws.Cells.Item(13, 14).Value = 1.119018180828921
$ws.Cells.Item(14, 3).Value = 1.112037276777103
$ws.Cells.Item(14, 4).Value = 1.114559537098774
$ws.Cells.Item(14, 5).Value = 1.111354557876756
$ws.Cells.Item(14, 6).Value = 1.124387784538375
$ws.Cells.Item(14, 10).Value = 1.118102857435521
$ws.Cells.Item(14, 11).Value = 1.11786134102969
$ws.Cells.Item(14, 12).Value = 1.114666521803918
$ws.Cells.Item(14, 13).Value = 1.127658826799744
$ws.Cells.Item(14, 14).Value = 1.11969069174127
$ws.Cells.Item(15, 3).Value = 1.112486816273379
$ws.Cells.Item(15, 4).Value = 1.114981972111544
$ws.Cells.Item(15, 5).Value = 1.111765264183701
$ws.Cells.Item(15, 6).Value = 1.124821334163839
$ws.Cells.Item(15, 10).Value = 1.118515932557199
$ws.Cells.Item(15, 11).Value = 1.118264505889452
$ws.Cells.Item(15, 12).Value = 1.115057928211827
$ws.Cells.Item(15, 13).Value = 1.128073273966431
$ws.Cells.Item(15, 14).Value = 1.120104353477001
$ws.Cells.Item(16, 3).Value = 1.11509257622649
$ws.Cells.Item(16, 4).Value = 1.117430389432967
$ws.Cells.Item(16, 5).Value = 1.114145576023771
$ws.Cells.Item(16, 6).Value = 1.127334324005294
$ws.Cells.Item(16, 10).Value = 1.120909796084228
$ws.Cells.Item(16, 11).Value = 1.120600807997515
$ws.Cells.Item(16, 12).Value = 1.117325943637882
$ws.Cells.Item(16, 13).Value = 1.130475127691016
$ws.Cells.Item(16, 14).Value = 1.122501616564819
$ws.Cells.Item(17, 3).Value = 1.116717830811897
$ws.Cells.Item(17, 4).Value = 1.118957302163238
$ws.Cells.Item(17, 5).Value = 1.115629904598439
$ws.Cells.Item(17, 6).Value = 1.128901640234791
$ws.Cells.Item(17, 10).Value = 1.122402420768214
$ws.Cells.Item(17, 11).Value = 1.122057428110147
$ws.Cells.Item(17, 12).Value = 1.118739858982833
$ws.Cells.Item(17, 13).Value = 1.131972768117045
$ws.Cells.Item(17, 14).Value = 1.123996360947062
$ws.Cells.Item(18, 3).Value = 1.117662552167203
$ws.Cells.Item(18, 4).Value = 1.119844785499603
$ws.Cells.Item(18, 5).Value = 1.116492596651917
$ws.Cells.Item(18, 6).Value = 1.129812655341856
$ws.Cells.Item(18, 10).Value = 1.123269879347642
$ws.Cells.Item(18, 11).Value = 1.122903921094696
$ws.Cells.Item(18, 12).Value = 1.119561487888792
$ws.Cells.Item(18, 13).Value = 1.132843153594857
$ws.Cells.Item(18, 14).Value = 1.124865051417172
$ws.Cells.Item(19, 3).Value = 1.117984131723657
$ws.Cells.Item(19, 4).Value = 1.120146868805906
$ws.Cells.Item(19, 5).Value = 1.116786234623619
$ws.Cells.Item(19, 6).Value = 1.130122756596414
$ws.Cells.Item(19, 10).Value = 1.123565130407842
$ws.Cells.Item(19, 11).Value = 1.123192029264542
$ws.Cells.Item(19, 12).Value = 1.119841125493468
$ws.Cells.Item(19, 13).Value = 1.133139402869376
$ws.Cells.Item(19, 14).Value = 1.125160721767744
$ws.Cells.Item(20, 3).Value = 1.116543795541707
$ws.Cells.Item(20, 4).Value = 1.118793805201207
$ws.Cells.Item(20, 5).Value = 1.115470971575041
$ws.Cells.Item(20, 6).Value = 1.128733812004516
$ws.Cells.Item(20, 10).Value = 1.122242605172773
$ws.Cells.Item(20, 11).Value = 1.121901471782345
$ws.Cells.Item(20, 12).Value = 1.118588479788396
$ws.Cells.Item(20, 13).Value = 1.131812414244664
$ws.Cells.Item(20, 14).Value = 1.123836318395144
$ws.Cells.Item(21, 3).Value = 1.111822230749608
$ws.Cells.Item(21, 4).Value = 1.114357452952329
$ws.Cells.Item(21, 5).Value = 1.111158082263438
$ws.Cells.Item(21, 6).Value = 1.124180386101842
$ws.Cells.Item(21, 10).Value = 1.117905245458829
$ws.Cells.Item(21, 11).Value = 1.117668467757862
$ws.Cells.Item(21, 12).Value = 1.114479271137417
$ws.Cells.Item(21, 13).Value = 1.127460559144061
$ws.Cells.Item(21, 14).Value = 1.119492799132904
$ws.Cells.Item(22, 3).Value = 1.108822501570186
$ws.Cells.Item(22, 4).Value = 1.111538262170439
$ws.Cells.Item(22, 5).Value = 1.108416985200761
$ws.Cells.Item(22, 6).Value = 1.121287238534618
$ws.Cells.Item(22, 10).Value = 1.115148083048277
$ws.Cells.Item(22, 11).Value = 1.114977269549086
$ws.Cells.Item(22, 12).Value = 1.111866352336084
$ws.Cells.Item(22, 13).Value = 1.124694295365074
$ws.Cells.Item(22, 14).Value = 1.116731721235479
$ws.Cells.Item(23, 3).Value = 1.11041574678678
$ws.Cells.Item(23, 4).Value = 1.1130356802934
$ws.Cells.Item(23, 5).Value = 1.109872957988025
$ws.Cells.Item(23, 6).Value = 1.122823897189452
$ws.Cells.Item(23, 10).Value = 1.116612638168604
$ws.Cells.Item(23, 11).Value = 1.116406820780208
$ws.Cells.Item(23, 12).Value = 1.113254362070876
$ws.Cells.Item(23, 13).Value = 1.126163674074597
$ws.Cells.Item(23, 14).Value = 1.118198356192063
$ws.Cells.Item(24, 3).Value = 1.116622444653475
$ws.Cells.Item(24, 4).Value = 1.11886769212567
$ws.Cells.Item(24, 5).Value = 1.115542796106206
$ws.Cells.Item(24, 6).Value = 1.128809656150222
$ws.Cells.Item(24, 10).Value = 1.12231482872387
$ws.Cells.Item(24, 11).Value = 1.121971951387869
$ws.Cells.Item(24, 12).Value = 1.11865689104603
$ws.Cells.Item(24, 13).Value = 1.131884881015935
$ws.Cells.Item(24, 14).Value = 1.123908644511968
$ws.Cells.Item(25, 3).Value = 1.123697798587814
$ws.Cells.Item(25, 4).Value = 1.125513024928253
$ws.Cells.Item(25, 5).Value = 1.12200177299882
$ws.Cells.Item(25, 6).Value = 1.135632050216936
$ws.Cells.Item(25, 10).Value = 1.128808519098931
$ws.Cells.Item(25, 11).Value = 1.128307958797856
$ws.Cells.Item(25, 12).Value = 1.124805950663291
$ws.Cells.Item(25, 13).Value = 1.138400689581405
$ws.Cells.Item(25, 14).Value = 1.130411556672198
